# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates the DAMSLTag (col I) and DialogAct (col J)
# values for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(14, 9).Value = 'sd'
$ws.Cells.Item(14, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(24, 9).Value = 'ba'
$ws.Cells.Item(24, 10).Value = 'Appreciation'
$ws.Cells.Item(27, 9).Value = 'b'
$ws.Cells.Item(27, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(84, 9).Value = 'b'
$ws.Cells.Item(84, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(94, 9).Value = 'ba'
$ws.Cells.Item(94, 10).Value = 'Appreciation'
$ws.Cells.Item(99, 9).Value = 'b'
$ws.Cells.Item(99, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(103, 9).Value = 'sd'
$ws.Cells.Item(103, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(122, 9).Value = 'ba'
$ws.Cells.Item(122, 10).Value = 'Appreciation'
$ws.Cells.Item(131, 9).Value = 'b'
$ws.Cells.Item(131, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(143, 9).Value = 'sd'
$ws.Cells.Item(143, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(145, 9).Value = '%'
$ws.Cells.Item(145, 10).Value = 'Uninterpretable'
$ws.Cells.Item(160, 9).Value = 'sv'
$ws.Cells.Item(160, 10).Value = 'Statement-opinion'
$ws.Cells.Item(161, 9).Value = 'sv'
$ws.Cells.Item(161, 10).Value = 'Statement-opinion'
$ws.Cells.Item(166, 9).Value = 'b'
$ws.Cells.Item(166, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(168, 9).Value = 'ba'
$ws.Cells.Item(168, 10).Value = 'Appreciation'
$ws.Cells.Item(172, 9).Value = 'sd'
$ws.Cells.Item(172, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(176, 9).Value = 'sd'
$ws.Cells.Item(176, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(180, 9).Value = 'sd'
$ws.Cells.Item(180, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(183, 9).Value = 'sd'
$ws.Cells.Item(183, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(199, 9).Value = 'b'
$ws.Cells.Item(199, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(203, 9).Value = 'b'
$ws.Cells.Item(203, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(205, 9).Value = 'sd'
$ws.Cells.Item(205, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(220, 9).Value = 'ba'
$ws.Cells.Item(220, 10).Value = 'Appreciation'
$ws.Cells.Item(233, 9).Value = 'sv'
$ws.Cells.Item(233, 10).Value = 'Statement-opinion'
$ws.Cells.Item(241, 9).Value = 'aa'
$ws.Cells.Item(241, 10).Value = 'Agree/Accept'
$ws.Cells.Item(242, 9).Value = 'b'
$ws.Cells.Item(242, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(249, 9).Value = 'sd'
$ws.Cells.Item(249, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(256, 9).Value = 'sv'
$ws.Cells.Item(256, 10).Value = 'Statement-opinion'
$ws.Cells.Item(260, 9).Value = 'sd'
$ws.Cells.Item(260, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(266, 9).Value = 'sd'
$ws.Cells.Item(266, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(267, 9).Value = 'b'
$ws.Cells.Item(267, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(280, 9).Value = 'aa'
$ws.Cells.Item(280, 10).Value = 'Agree/Accept'
$ws.Cells.Item(290, 9).Value = 'sd'
$ws.Cells.Item(290, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(326, 9).Value = 'ba'
$ws.Cells.Item(326, 10).Value = 'Appreciation'
$ws.Cells.Item(340, 9).Value = 'sv'
$ws.Cells.Item(340, 10).Value = 'Statement-opinion'
$ws.Cells.Item(347, 9).Value = 'b'
$ws.Cells.Item(347, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(356, 9).Value = 'b'
$ws.Cells.Item(356, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(357, 9).Value = 'sv'
$ws.Cells.Item(357, 10).Value = 'Statement-opinion'
$ws.Cells.Item(380, 9).Value = 'sd'
$ws.Cells.Item(380, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(384, 9).Value = 'sd'
$ws.Cells.Item(384, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(415, 9).Value = 'sv'
$ws.Cells.Item(415, 10).Value = 'Statement-opinion'
$ws.Cells.Item(418, 9).Value = 'sd'
$ws.Cells.Item(418, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(420, 9).Value = 'sv'
$ws.Cells.Item(420, 10).Value = 'Statement-opinion'
$ws.Cells.Item(425, 9).Value = 'b'
$ws.Cells.Item(425, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(429, 9).Value = 'b'
$ws.Cells.Item(429, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(443, 9).Value = 'sd'
$ws.Cells.Item(443, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(453, 9).Value = 'b'
$ws.Cells.Item(453, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(460, 9).Value = 'sd'
$ws.Cells.Item(460, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(464, 9).Value = 'sd'
$ws.Cells.Item(464, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(470, 9).Value = 'sd'
$ws.Cells.Item(470, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(473, 9).Value = 'b'
$ws.Cells.Item(473, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(483, 9).Value = 'sd'
$ws.Cells.Item(483, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(504, 9).Value = 'sd'
$ws.Cells.Item(504, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(505, 9).Value = 'ba'
$ws.Cells.Item(505, 10).Value = 'Appreciation'
$ws.Cells.Item(521, 9).Value = 'sv'
$ws.Cells.Item(521, 10).Value = 'Statement-opinion'
$ws.Cells.Item(523, 9).Value = 'aa'
$ws.Cells.Item(523, 10).Value = 'Agree/Accept'
$ws.Cells.Item(547, 9).Value = 'aa'
$ws.Cells.Item(547, 10).Value = 'Agree/Accept'